$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-05 Monday" "2024-02-06 Tuesday"

Replace-Text "25×38=950" "84×86=7224"
Replace-Text "24×21=504" "94×78=7332"
Replace-Text "87×72=6264" "84×86=7224"
Replace-Text "22×51=1122" "74×43=3182"
Replace-Text "94×71=6674" "53×45=2385"

Replace-Text "88×74=6512" "57×84=4788"
Replace-Text "92×74=6808" "91×88=8008"
Replace-Text "27×58=1566" "12×78=936"
Replace-Text "89×17=1513" "65×88=5720"
Replace-Text "51×24=1224" "91×21=1911"

Replace-Text "67×85=5695" "85×28=2380"
Replace-Text "47×94=4418" "82×96=7872"
Replace-Text "37×82=3034" "54×54=2916"
Replace-Text "26×44=1144" "66×33=2178"
Replace-Text "22×35=770" "17×69=1173"

Replace-Text "56×63=3528" "57×66=3762"
Replace-Text "39×19=741" "37×55=2035"
Replace-Text "53×89=4717" "75×91=6825"
Replace-Text "92×93=8556" "91×68=6188"
Replace-Text "46×87=4002" "40×35=1400"

Replace-Text "58×96=5568" "75×84=6300"
Replace-Text "12×50=600" "72×85=6120"
Replace-Text "41×57=2337" "72×84=6048"
Replace-Text "66×95=6270" "30×43=1290"
Replace-Text "78×68=5304" "65×47=3055"
